$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update the input value (hyper extended knee fix): C2 changes from 52 to 122.
# Dependent formulas (D2, G2, H2, I2, J2) recalculate automatically.
$ws.Range("C2").Value = 122

$excel.CalculateFullRebuild()

# Make sure the chart that plots column D ('Ark1'!$D$2:$D$9) picks up the
# refreshed value for its cached data point.
$co = $ws.ChartObjects().Item(1)
$co.Chart.Refresh()

# Update the active selection to reflect the saved view state (C2 selected).
$ws.Activate()
$ws.Range("C2").Select()
